# Update the "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect freshly generated data (gh-pages output refresh).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (row on sheet -> new value for column F)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 571
$ws1.Range("F4").Value = 1243
$ws1.Range("F5").Value = 1066
$ws1.Range("F6").Value = 14031
$ws1.Range("F7").Value = 15373
$ws1.Range("F9").Value = 34
$ws1.Range("F20").Value = 1188
$ws1.Range("F23").Value = 5979
$ws1.Range("F25").Value = 1078
$ws1.Range("F26").Value = 5519
$ws1.Range("F27").Value = 67
$ws1.Range("F29").Value = 98
$ws1.Range("F30").Value = 444

# Sheet "全部类型" (same underlying rows, shifted down by one)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 571
$ws4.Range("F5").Value = 1243
$ws4.Range("F6").Value = 1066
$ws4.Range("F7").Value = 14031
$ws4.Range("F8").Value = 15373
$ws4.Range("F10").Value = 34
$ws4.Range("F21").Value = 1188
$ws4.Range("F25").Value = 5979
$ws4.Range("F27").Value = 1078
$ws4.Range("F28").Value = 5519
$ws4.Range("F29").Value = 67
$ws4.Range("F31").Value = 98
$ws4.Range("F32").Value = 444
